# Update the data table:
#  - Replace the sample rows (A2:J19) with the refreshed dataset values
#  - Drop the trailing rows (20-22) that are no longer part of the sample
#    (they were quantile-filtered outliers / extra rows removed upstream)
#  - Shrinks the sheet dimension automatically to A1:J19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 18,10
$data[0,0] = 129
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 13
$data[0,5] = 1
$data[0,6] = 3
$data[0,7] = 40
$data[0,8] = 18
$data[0,9] = 0
$data[1,0] = 175
$data[1,1] = 0
$data[1,2] = 1
$data[1,3] = 0
$data[1,4] = 5
$data[1,5] = 1
$data[1,6] = 6
$data[1,7] = 40
$data[1,8] = 21
$data[1,9] = 0
$data[2,0] = 202
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 1
$data[2,4] = 2
$data[2,5] = 1
$data[2,6] = 7
$data[2,7] = 49
$data[2,8] = 5
$data[2,9] = 0
$data[3,0] = 19
$data[3,1] = 1
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 13
$data[3,5] = 0
$data[3,6] = 6
$data[3,7] = 39
$data[3,8] = 5
$data[3,9] = 0
$data[4,0] = 28
$data[4,1] = 0
$data[4,2] = 1
$data[4,3] = 0
$data[4,4] = 10
$data[4,5] = 1
$data[4,6] = 12
$data[4,7] = 57
$data[4,8] = 16
$data[4,9] = 0
$data[5,0] = 260
$data[5,1] = 0
$data[5,2] = 0
$data[5,3] = 1
$data[5,4] = 13
$data[5,5] = 0
$data[5,6] = 4
$data[5,7] = 36
$data[5,8] = 20
$data[5,9] = 0
$data[6,0] = 206
$data[6,1] = 0
$data[6,2] = 1
$data[6,3] = 1
$data[6,4] = 2
$data[6,5] = 0
$data[6,6] = 1
$data[6,7] = 30
$data[6,8] = 12
$data[6,9] = 1
$data[7,0] = 272
$data[7,1] = 0
$data[7,2] = 1
$data[7,3] = 0
$data[7,4] = 12
$data[7,5] = 1
$data[7,6] = 11
$data[7,7] = 55
$data[7,8] = 4
$data[7,9] = 1
$data[8,0] = 87
$data[8,1] = 0
$data[8,2] = 0
$data[8,3] = 0
$data[8,4] = 13
$data[8,5] = 0
$data[8,6] = 10
$data[8,7] = 53
$data[8,8] = 13
$data[8,9] = 0
$data[9,0] = 0
$data[9,1] = 0
$data[9,2] = 1
$data[9,3] = 1
$data[9,4] = 1.5
$data[9,5] = 0
$data[9,6] = 8
$data[9,7] = 45
$data[9,8] = 5
$data[9,9] = 0
$data[10,0] = 64
$data[10,1] = 0
$data[10,2] = 1
$data[10,3] = 0
$data[10,4] = 6
$data[10,5] = 0
$data[10,6] = 2
$data[10,7] = 46
$data[10,8] = 2
$data[10,9] = 0
$data[11,0] = 298
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 5
$data[11,5] = 0
$data[11,6] = 5
$data[11,7] = 37
$data[11,8] = 13
$data[11,9] = 0
$data[12,0] = 161
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 0
$data[12,4] = 1
$data[12,5] = 1
$data[12,6] = 5
$data[12,7] = 36
$data[12,8] = 8
$data[12,9] = 0
$data[13,0] = 45
$data[13,1] = 0
$data[13,2] = 1
$data[13,3] = 0
$data[13,4] = 5
$data[13,5] = 1
$data[13,6] = 7
$data[13,7] = 53
$data[13,8] = 2
$data[13,9] = 0
$data[14,0] = 22
$data[14,1] = 0
$data[14,2] = 1
$data[14,3] = 1
$data[14,4] = 8
$data[14,5] = 0
$data[14,6] = 8
$data[14,7] = 45
$data[14,8] = 7
$data[14,9] = 0
$data[15,0] = 18
$data[15,1] = 0
$data[15,2] = 1
$data[15,3] = 0
$data[15,4] = 12
$data[15,5] = 1
$data[15,6] = 11
$data[15,7] = 51
$data[15,8] = 7
$data[15,9] = 0
$data[16,0] = 262
$data[16,1] = 0
$data[16,2] = 1
$data[16,3] = 0
$data[16,4] = 4
$data[16,5] = 1
$data[16,6] = 10
$data[16,7] = 56
$data[16,8] = 10
$data[16,9] = 0
$data[17,0] = 54
$data[17,1] = 0
$data[17,2] = 1
$data[17,3] = 1
$data[17,4] = 2
$data[17,5] = 0
$data[17,6] = 11
$data[17,7] = 54
$data[17,8] = 4
$data[17,9] = 0

$ws.Range("A2:J19").Value = $data

# Remove the now-unused trailing rows so the used range / dimension shrinks to J19
$ws.Rows("20:22").Delete()
